$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new record (row 9) to the log: check-out entry for Ahmad Sharim
# recorded at 02/05/2025 02:09:38 PM.
$row = 9

$ws.Cells.Item($row, 1).Value = ""
$ws.Cells.Item($row, 2).Value = "أحمد شريم"

# Column C ("الكمية") holds values that look numeric but are stored as text
# throughout the sheet (see C2:C8). Enter it as a formula that yields the
# text "333" and then paste the computed result back as a plain value so it
# lands in the cell as text (not a number) without attaching any new/explicit
# number-format style to the cell.
$ws.Cells.Item($row, 3).Formula = "=""333"""
$ws.Cells.Item($row, 3).Copy() | Out-Null
$ws.Cells.Item($row, 3).PasteSpecial(-4163) | Out-Null

$ws.Cells.Item($row, 4).Value = "النصر"
$ws.Cells.Item($row, 5).Value = "الرحلة 3"
$ws.Cells.Item($row, 6).Value = "C5"
$ws.Cells.Item($row, 7).Value = "WCK"
$ws.Cells.Item($row, 8).Value = "٠٢‏/٠٥‏/٢٠٢٥ ٠٢:٠٩:٣٨ م"
